{"js": "// \"added title to mttw\" \u2014 append a title to the \"Move to the Web Team\"\n// experience heading so it reads \"Move to the Web Team \u2013 Full Stack Web\n// Developer\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the exact paragraph whose text is \"Move to the Web Team\" (the\n// job-title line for that experience entry \u2014 distinct from the later\n// \"Data Move to the Web Team \u2013 ...\" paragraph).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Move to the Web Team\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Move to the Web Team\"');\n}\n\n// Find the trailing \"Team\" word inside that paragraph so the new title text\n// is inserted right after it, inheriting its character formatting\n// (9pt / Calibri complex-script font) automatically.\nconst found = target.search(\"Team\", { matchCase: true, matchWholeWord: true });\nawait context.sync();\n\nif (found.items.length === 0) {\n  throw new Error('Could not find \"Team\" inside the target paragraph');\n}\n\nconst teamRange = found.items[found.items.length - 1];\nteamRange.insertText(\" \\u2013 Full Stack Web Developer\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# \"added title to mttw\" \u2014 append a title to the \"Move to the Web Team\"\n# experience heading so it reads \"Move to the Web Team \u2013 Full Stack Web\n# Developer\".\n\n$d = $word.ActiveDocument\n\n# Locate the exact paragraph whose text is \"Move to the Web Team\" (the\n# job-title line for that experience entry \u2014 distinct from the later\n# \"Data Move to the Web Team - ...\" paragraph).\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text.Trim() -eq \"Move to the Web Team\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find paragraph 'Move to the Web Team'\"\n}\n\n$enDash = [char]0x2013\n$replacementText = \"Team \" + $enDash + \" Full Stack Web Developer\"\n\n# Find the trailing \"Team\" word inside that paragraph and replace it with\n# \"Team - Full Stack Web Developer\" in one go, so the appended text\n# inherits the matched run's character formatting (9pt / Calibri\n# complex-script font) automatically.\n$find = $target.Range.Find\n$find.ClearFormatting()\n$find.Execute(\"Team\", $true, $true, $false, $false, $false, $true, 1, $false, $replacementText, 2) | Out-Null\n"}
